$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row right above the trailing blank-formatted row (old row 10).
# Excel's native row-insert shifts row 10's content/format down to row 11
# intact, and auto-fills the new row 10 with row 9's formatting.
$ws.Rows(10).Insert()

# Grow the table definition to include the newly inserted row.
$lo.Resize($ws.Range("A1:C11"))

# Row 9 used to hold the now-removed "Regex" entry; turn it into the new
# AttachmentDirectory parameter row (match formatting of the other data rows,
# using row 6 - a wrapped/tall description row - as the template).
$ws.Range("A6:C6").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows(9).RowHeight = $ws.Rows(6).RowHeight

$ws.Range("A9").Value2 = "AttachmentDirectory"
$ws.Range("B9").Value2 = "{0}\Desktop\Attachments"
$ws.Range("C9").Value2 = "The file path for storing email attachments"

$excel.CutCopyMode = $false

# Update the selected cell as recorded in the saved workbook.
$ws.Range("C16").Select()
